$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# Update Date value
$wsMetadata.Range("B8").Value = "2025-04-30T13:43:05+00:00"

# Update Description value
$wsMetadata.Range("B12").Value = "Model logique d'un auteur"

# Update Max / Base Max for Author.specialty row (row 6) from "*" to "1".
# Plain assignment (.Value = "1") would store the cell as a number, but the
# target keeps it as a shared string. Copy/PasteSpecial(Values) from an
# existing cell that already holds the text "1" (same column, another row)
# preserves the text type without touching the cell's style.
$wsElements.Range("G7").Copy()
$wsElements.Range("G6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$wsElements.Range("AH7").Copy()
$wsElements.Range("AH6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$excel.CutCopyMode = $false
